$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (C) and montant_total (D) for rows with revised 2020-09-02 data
# Keep values stored as text to match the source data format (text-formatted numeric strings)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "193"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "454016.00"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "44"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "88000.00"

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "102"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "283752.38"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "431"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1376161.23"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "156"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "675202.10"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "53"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "128800.00"

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "49"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "131669.13"

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "140"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "444965.56"

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "62"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "242297.00"

$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "589"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1941776.10"

$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "247"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1303520.11"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "174"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "480779.00"

$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "107"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "321135.17"

$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "612"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "2226761.21"

$ws.Range("C53").NumberFormat = "@"
$ws.Range("C53").Value = "270"
$ws.Range("D53").NumberFormat = "@"
$ws.Range("D53").Value = "1222878.76"

$ws.Range("C54").NumberFormat = "@"
$ws.Range("C54").Value = "93"
$ws.Range("D54").NumberFormat = "@"
$ws.Range("D54").Value = "558274.23"

$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").Value = "717"
$ws.Range("D57").NumberFormat = "@"
$ws.Range("D57").Value = "1850818.62"

$ws.Range("C64").NumberFormat = "@"
$ws.Range("C64").Value = "24"
$ws.Range("D64").NumberFormat = "@"
$ws.Range("D64").Value = "60261.00"

$ws.Range("C65").NumberFormat = "@"
$ws.Range("C65").Value = "115"
$ws.Range("D65").NumberFormat = "@"
$ws.Range("D65").Value = "281349.69"

$ws.Range("C66").NumberFormat = "@"
$ws.Range("C66").Value = "62"
$ws.Range("D66").NumberFormat = "@"
$ws.Range("D66").Value = "203035.00"

$ws.Range("C67").NumberFormat = "@"
$ws.Range("C67").Value = "20"
$ws.Range("D67").NumberFormat = "@"
$ws.Range("D67").Value = "88876.00"

$ws.Range("C68").NumberFormat = "@"
$ws.Range("C68").Value = "4"
$ws.Range("D68").NumberFormat = "@"
$ws.Range("D68").Value = "12500.00"

$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "908"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "2933833.26"

$ws.Range("C85").NumberFormat = "@"
$ws.Range("C85").Value = "338"
$ws.Range("D85").NumberFormat = "@"
$ws.Range("D85").Value = "1373803.79"

$ws.Range("C86").NumberFormat = "@"
$ws.Range("C86").Value = "120"
$ws.Range("D86").NumberFormat = "@"
$ws.Range("D86").Value = "598484.52"
